$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = "328.40";       E = "-0.94%" },
    @{ Row = 3;  D = "43.60";        E = "4.72%" },
    @{ Row = 4;  D = "5.563";        E = "-2.09%" },
    @{ Row = 5;  D = "0.08187";      E = "-2.28%" },
    @{ Row = 6;  D = "8.754";        E = "-0.48%" },
    @{ Row = 7;  D = "4.359";        E = "-3.64%" },
    @{ Row = 8;  D = "1.876";        E = "-7.02%" },
    @{ Row = 9;  D = "2.801";        E = "-3.12%" },
    @{ Row = 10; D = "0.9441";       E = "1.94%" },
    @{ Row = 11; D = "0.1202";       E = "-7.01%" },
    @{ Row = 12; D = "0.1905";       E = "-3.37%" },
    @{ Row = 13; D = "0.09751";      E = "3.54%" },
    @{ Row = 14; D = "0.04161";      E = "6.98%" },
    @{ Row = 15; E = "0.87%" },
    @{ Row = 16; D = "0.001295";     E = "-1.01%" },
    @{ Row = 17; D = "0.006062";     E = "-0.61%" },
    @{ Row = 18; D = "3.537";        E = "2.96%" },
    @{ Row = 20; E = "7.21%" },
    @{ Row = 21; E = "-0.23%" },
    @{ Row = 22; D = "0.2503";       E = "-0.30%" },
    @{ Row = 23; D = "0.04394";      E = "-0.47%" },
    @{ Row = 24; D = "0.001242";     E = "-1.23%" },
    @{ Row = 25; D = "0.004293";     E = "-2.58%" },
    @{ Row = 26; D = "0.0001238";    E = "3.12%" },
    @{ Row = 27; D = "0.0004016";    E = "31.87%" },
    @{ Row = 39; D = "0.02724";      E = "-3.35%" },
    @{ Row = 40; D = "0.05691";      E = "2.89%" },
    @{ Row = 41; D = "0.007886";     E = "1.20%" },
    @{ Row = 42; D = "0.009773";     E = "4.85%" },
    @{ Row = 43; D = "0.1413";       E = "-1.70%" },
    @{ Row = 44; D = "0.002109";     E = "-2.34%" },
    @{ Row = 45; D = "0.01001";      E = "-9.59%" },
    @{ Row = 46; D = "0.00007334";   E = "4.51%" },
    @{ Row = 47; D = "0.00000000755"; E = "0.58%" },
    @{ Row = 48; D = "0.003447";     E = "-2.16%" },
    @{ Row = 49; D = "0.002285";     E = "0.24%" },
    @{ Row = 50; D = "0.00002114";   E = "0.58%" },
    @{ Row = 51; D = "0.0002013";    E = "0.58%" }
)

foreach ($u in $updates) {
    if ($u.ContainsKey("D")) {
        $cell = $ws.Range("D$($u.Row)")
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
    }
    if ($u.ContainsKey("E")) {
        $cell = $ws.Range("E$($u.Row)")
        $cell.NumberFormat = "@"
        $cell.Value = $u.E
    }
}
